# Auto-committed on 2023/02/03 週五 17:21:57.66
# GenTable總表 refresh: update a handful of "last modified" timestamps for
# existing tables, and insert the newly generated table TxAmlRatingAppl
# (alphabetically, right after TxAmlRating) into the XX-系統 section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh "最後修改時間" (last modified time) for tables that changed.
# ---------------------------------------------------------------------
$timestampUpdates = @{
    "CustDataCtrl" = "2023年02月02日 19:28:47"
    "FacMain"      = "2023年01月30日 10:58:58"
    "InsuRenew"    = "2023年02月02日 17:20:40"
    "NegFinShare"  = "2023年02月01日 17:57:44"
    "AcReceivable" = "2023年01月31日 14:07:34"
    "CdReport"     = "2023年01月30日 15:18:39"
    "SystemParas"  = "2023年02月01日 10:10:56"
    "DailyTav"     = "2023年02月03日 10:56:15"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $tableName = $ws.Cells.Item($r, 2).Value2
    if ($timestampUpdates.ContainsKey($tableName)) {
        $ws.Cells.Item($r, 5).Value = $timestampUpdates[$tableName]
    }
}

# ---------------------------------------------------------------------
# 2) Insert the new TxAmlRatingAppl row, right after TxAmlRating (row 339)
#    so the XX-系統 section stays sorted alphabetically. Every row from
#    340 downward shifts down by one.
# ---------------------------------------------------------------------
$ws.Rows("340").Insert()

$ws.Cells.Item(340, 1).Value = "XX-系統"
$ws.Cells.Item(340, 2).Value = "TxAmlRatingAppl"
$ws.Cells.Item(340, 3).Value = "Eloan評級案件申請留存檔"
$ws.Cells.Item(340, 4).Formula = '=HYPERLINK("[\\192.168.10.16\St1Share(NAS)\SKL\DB\GenTables\XX-系統\TxAmlRatingAppl.xlsx]DBD!A1", "連結")'
$ws.Cells.Item(340, 5).Value = "2023年02月03日 09:11:36"
